$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.260.98'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.689.71'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = "'219.30"
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").Value = "'0.5243"
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +2.05%  '
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").Value = "'22.06"
$ws.Range("E10").Value = '  +2.39%  '
$ws.Range("D11").Value = "'0.07476"
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").Value = '1.705.58'
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").Value = "'4.554"
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("D14").Value = "'0.5859"
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("D15").Value = "'0.000008554"
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").Value = "'64.67"
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '26.317.33'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").Value = "'4.970"
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").Value = "'190.61"
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").Value = "'6.236"
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").Value = "'1.007"
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = "'145.32"
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("D25").Value = "'7.681"
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").Value = "'0.1245"
$ws.Range("E26").Value = '  +5.85%  '
$ws.Range("D27").Value = "'15.89"
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = "'0.06691"
$ws.Range("E28").Value = '  +14.25%  '
$ws.Range("D29").Value = "'1.351"
$ws.Range("E29").Value = '  +4.34%  '
$ws.Range("D30").Value = "'1.332"
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("D31").Value = "'3.603"
$ws.Range("E31").Value = '  +2.74%  '
$ws.Range("D32").Value = "'3.552"
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("D34").Value = "'1.028"
$ws.Range("E34").Value = '  +1.88%  '
$ws.Range("D35").Value = "'0.6207"
$ws.Range("E35").Value = '  +3.38%  '
$ws.Range("D36").Value = "'2.389"
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("E37").Value = '  +2.39%  '
$ws.Range("D38").Value = "'6.283"
$ws.Range("E38").Value = '  +5.53%  '
$ws.Range("D39").Value = "'0.01619"
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("D40").Value = '1.103.38'
$ws.Range("E40").Value = '  +1.46%  '
$ws.Range("D41").Value = "'0.8802"
$ws.Range("E41").Value = '  +2.30%  '
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("D43").Value = "'100.96"
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").Value = '1.837.39'
$ws.Range("E44").Value = '  +1.11%  '
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("D46").Value = "'56.88"
$ws.Range("E46").Value = '  +1.65%  '
# Row 47/48: EnergySwap and Frax swap positions, with updated price/volume data
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = "'1.009"
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'8.158"
$ws.Range("E48").Value = '  +1.28%  '
$ws.Range("D49").Value = "'0.05264"
$ws.Range("E49").Value = '  +1.72%  '
$ws.Range("D50").Value = "'0.4298"
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("E51").Value = '  +1.84%  '
